{"js": "// Append a parenthetical red-colored note to the first paragraph of the\n// document (\"This is a Microsoft word document.\"), turning it into:\n//   \"This is a Microsoft word document.  (This is a change \u2013 Version for main branch)\"\n// with the trailing note rendered in red (FF0000), split across three runs\n// the way an interactive typing/paste session would naturally produce.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// 1. Add two trailing spaces to the existing (black) sentence, keeping the\n//    original run's formatting untouched.\nfirstParagraph.insertText(\"  \", Word.InsertLocation.end);\n\n// 2. Insert the red parenthetical note. It is inserted as three separate\n//    runs (matching how the change was authored) but each chunk gets the\n//    same red font color applied right after insertion.\nconst redHex = \"#FF0000\";\n\nconst redRun1 = firstParagraph.insertText(\n  \"(This is a change \\u2013 Ve\",\n  Word.InsertLocation.end\n);\nredRun1.font.color = redHex;\n\nconst redRun2 = firstParagraph.insertText(\n  \"rsion for main branch\",\n  Word.InsertLocation.end\n);\nredRun2.font.color = redHex;\n\nconst redRun3 = firstParagraph.insertText(\")\", Word.InsertLocation.end);\nredRun3.font.color = redHex;\n\nawait context.sync();\n", "ps1": "# Append a parenthetical red-colored note to the first paragraph of the\n# document (\"This is a Microsoft word document.\"), turning it into:\n#   \"This is a Microsoft word document.  (This is a change \u2013 Version for main branch)\"\n# with the trailing note rendered in red (RGB 255,0,0 -> WdColor 255) and\n# split across three runs, matching how the change was authored.\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# Range of the first paragraph, excluding the trailing paragraph mark.\n$r = $p.Range\n$r.MoveEnd(1, -1)   # wdCharacter = 1\n\n# 1. Add two trailing spaces after the existing (black) sentence, keeping\n#    the original formatting.\n$r.InsertAfter(\"  \")\n$r.Collapse(0)      # wdCollapseEnd = 0\n\n$wdColorRed = 255    # WdColor enum value for pure red (BGR 0x0000FF)\n\n# 2. Insert the red parenthetical note as three separate runs.\n$r2 = $d.Range($r.End, $r.End)\n$r2.InsertAfter(\"(This is a change \u2013 Ve\")\n$r2.Font.Color = $wdColorRed\n$r2.Collapse(0)\n\n$r3 = $d.Range($r2.End, $r2.End)\n$r3.InsertAfter(\"rsion for main branch\")\n$r3.Font.Color = $wdColorRed\n$r3.Collapse(0)\n\n$r4 = $d.Range($r3.End, $r3.End)\n$r4.InsertAfter(\")\")\n$r4.Font.Color = $wdColorRed\n"}
